$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'298.58"
$ws.Range("E2").Value = "'-1.10%"
$ws.Range("D3").Value = "'31.76"
$ws.Range("E3").Value = "'0.79%"
$ws.Range("D4").Value = "'5.087"
$ws.Range("E4").Value = "'-1.25%"
$ws.Range("D5").Value = "'0.08153"
$ws.Range("E5").Value = "'10.57%"
$ws.Range("D6").Value = "'2.568"
$ws.Range("E6").Value = "'-1.27%"
$ws.Range("D7").Value = "'7.771"
$ws.Range("E7").Value = "'-1.76%"
$ws.Range("D8").Value = "'3.846"
$ws.Range("E8").Value = "'2.19%"
$ws.Range("D9").Value = "'0.9317"
$ws.Range("E9").Value = "'1.30%"
$ws.Range("D10").Value = "'0.1760"
$ws.Range("E10").Value = "'2.21%"
$ws.Range("D11").Value = "'0.07516"
$ws.Range("E11").Value = "'0.70%"
$ws.Range("D12").Value = "'0.09096"
$ws.Range("E12").Value = "'11.88%"
$ws.Range("D13").Value = "'0.02999"
$ws.Range("E13").Value = "'-1.36%"
$ws.Range("E14").Value = "'0.78%"
$ws.Range("D15").Value = "'0.001510"
$ws.Range("E15").Value = "'0.71%"
$ws.Range("D16").Value = "'0.005935"
$ws.Range("E16").Value = "'-3.64%"
$ws.Range("D17").Value = "'3.574"
$ws.Range("E17").Value = "'3.63%"
$ws.Range("D18").Value = "'2.255"
$ws.Range("E18").Value = "'1.15%"
$ws.Range("E19").Value = "'-0.53%"
$ws.Range("D20").Value = "'0.1327"
$ws.Range("E20").Value = "'-0.01%"
$ws.Range("D21").Value = "'3.951"
$ws.Range("E21").Value = "'-15.22%"
$ws.Range("D22").Value = "'0.1699"
$ws.Range("E22").Value = "'7.04%"
$ws.Range("D23").Value = "'0.04602"
$ws.Range("E23").Value = "'-0.96%"
$ws.Range("D24").Value = "'0.001241"
$ws.Range("E24").Value = "'0.92%"
$ws.Range("D25").Value = "'0.004464"
$ws.Range("E25").Value = "'-0.10%"
$ws.Range("D26").Value = "'0.0001199"
$ws.Range("E26").Value = "'-7.89%"
$ws.Range("D27").Value = "'0.0003410"
$ws.Range("E27").Value = "'82.20%"
$ws.Range("D39").Value = "'0.01769"
$ws.Range("E39").Value = "'3.24%"
$ws.Range("D40").Value = "'0.04539"
$ws.Range("E40").Value = "'0.35%"
$ws.Range("D41").Value = "'0.006988"
$ws.Range("E41").Value = "'-1.75%"
$ws.Range("D42").Value = "'0.1355"
$ws.Range("E42").Value = "'0.86%"
$ws.Range("D43").Value = "'0.002208"
$ws.Range("E43").Value = "'-1.11%"
$ws.Range("D44").Value = "'0.009941"
$ws.Range("E44").Value = "'-9.04%"
$ws.Range("D45").Value = "'0.00006437"
$ws.Range("E45").Value = "'2.04%"
$ws.Range("E46").Value = "'-0.09%"
$ws.Range("D47").Value = "'0.008739"
$ws.Range("E47").Value = "'-12.94%"
$ws.Range("E48").Value = "'11.10%"
$ws.Range("D49").Value = "'0.00002098"
$ws.Range("E49").Value = "'-0.09%"
$ws.Range("D50").Value = "'0.0001998"
$ws.Range("E50").Value = "'-0.02%"
